$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as text (matching the original formatting) instead
# of silently converting to a number.
$textCells = @("D5","D10","D11","D15","D17","D20","D22","D24","D25","D26","D28","D29","D31","D35","D36","D38","D40","D41","D42","D43","D46","D47","D49","D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.061.03"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.647.65"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  +0.93%  "
$ws.Range("D5").Value = "216.90"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("D10").Value = "19.71"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").Value = "1.875.10"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "1.653.07"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "0.547"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "63.25"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "26.038.45"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "193.45"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "9.97"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "0.133"
$ws.Range("E24").Value = "  +7.06%  "
$ws.Range("D25").Value = "1.81"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "144.43"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").Value = "6.95"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "0.0500"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("D36").Value = "0.909"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").Value = "1.133.62"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "0.544"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "0.0157"
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").Value = "5.52"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "99.72"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "0.799"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "1.783.56"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("D46").Value = "56.81"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "0.0535"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Value = "7.74"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "0.417"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("E51").Value = "  +0.15%  "

